$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the B column values (B1:B8) with the new decimal numbers.
$newValues = @(66661, 66782, 66663, 66666, 66511, 66444, 66757, 66469)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}

# Update the selection to cover C1:C8 as in the final sheetView.
$ws.Range("C1:C8").Select()
